$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "59÷9=6, 5"
$t.Cell(1, 2).Range.Text = "91÷4=22, 3"
$t.Cell(1, 3).Range.Text = "23÷4=5, 3"
$t.Cell(1, 4).Range.Text = "59÷3=19, 2"
$t.Cell(1, 5).Range.Text = "41÷6=6, 5"
$t.Cell(5, 1).Range.Text = "44÷2=22, 0"
$t.Cell(5, 2).Range.Text = "56÷3=18, 2"
$t.Cell(5, 3).Range.Text = "87÷9=9, 6"
$t.Cell(5, 4).Range.Text = "87÷4=21, 3"
$t.Cell(5, 5).Range.Text = "60÷8=7, 4"
$t.Cell(9, 1).Range.Text = "96÷5=19, 1"
$t.Cell(9, 2).Range.Text = "81÷9=9, 0"
$t.Cell(9, 3).Range.Text = "85÷3=28, 1"
$t.Cell(9, 4).Range.Text = "98÷2=49, 0"
$t.Cell(9, 5).Range.Text = "78÷4=19, 2"
$t.Cell(13, 1).Range.Text = "61÷5=12, 1"
$t.Cell(13, 2).Range.Text = "54÷7=7, 5"
$t.Cell(13, 3).Range.Text = "60÷4=15, 0"
$t.Cell(13, 4).Range.Text = "36÷3=12, 0"
$t.Cell(13, 5).Range.Text = "37÷2=18, 1"
$t.Cell(17, 1).Range.Text = "47÷7=6, 5"
$t.Cell(17, 2).Range.Text = "69÷5=13, 4"
$t.Cell(17, 3).Range.Text = "10÷7=1, 3"
$t.Cell(17, 4).Range.Text = "40÷2=20, 0"
$t.Cell(17, 5).Range.Text = "33÷3=11, 0"
